# Update countries & provincias Spain
# Applies the data refresh captured by the commit diff:
#  - Updates the "Datos actualizados" timestamp string
#  - Updates case numbers for several countries
#  - Armenia overtakes Kirguistan (Armenia's row now has fresh numbers,
#    Kirguistan keeps its previous numbers but drops one rank)
#  - Hungria overtakes Zimbabue (same pattern)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "last refreshed" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 2 de Septiembre de 2020 a las 09:13"

# Estados Unidos (row 4) - numbers refreshed
$ws.Range("B4").Value = 6257938
$ws.Range("C4").Value = 367
$ws.Range("D4").Value = 3497431
$ws.Range("E4").Value = 2571605
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 188902

# Ucrania (row 28) - numbers refreshed
$ws.Range("B28").Value = 125789
$ws.Range("C28").Value = 2495
$ws.Range("D28").Value = 58817
$ws.Range("E28").Value = 64316
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 51
$ws.Range("H28").Value = 2656

# Armenia moves up to row 59 with fresh numbers
$ws.Range("A59").Value = "Armenia"
$ws.Range("B59").Value = 44075
$ws.Range("C59").Value = 197
$ws.Range("D59").Value = 38631
$ws.Range("E59").Value = 4560
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = 884

# Kirguistan drops to row 60, keeping its previous (unchanged) numbers
$ws.Range("A60").Value = "Kirguistan"
$ws.Range("B60").Value = 44036
$ws.Range("C60").Value = 78
$ws.Range("D60").Value = 38895
$ws.Range("E60").Value = 4082
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 1059

# El Salvador (row 73) - numbers refreshed
$ws.Range("B73").Value = 25904
$ws.Range("C73").Value = 84
$ws.Range("D73").Value = 14879
$ws.Range("E73").Value = 10294
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 7
$ws.Range("H73").Value = 731

# Hungria moves up to row 106 with fresh numbers
$ws.Range("A106").Value = "Hungria"
$ws.Range("B106").Value = 6622
$ws.Range("C106").Value = 365
$ws.Range("D106").Value = 3903
$ws.Range("E106").Value = 2100
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 3
$ws.Range("H106").Value = 619

# Zimbabue drops to row 107, keeping its previous (unchanged) numbers
$ws.Range("A107").Value = "Zimbabue"
$ws.Range("B107").Value = 6559
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 5241
$ws.Range("E107").Value = 1115
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 203

# Georgia (row 152) - numbers refreshed
$ws.Range("B152").Value = 1548
$ws.Range("C152").Value = 38
$ws.Range("D152").Value = 1270
$ws.Range("E152").Value = 259
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 19
